$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse the 2x2 reaction-term matrix (a(t),b(t),c(t),d(t) with
# a(t)*u+b(t)*v / c(t)*u+d(t)*v) down to a plain "u+v" sum, used for both
# the Hydon and DBH error-search matrix definitions.

# D2 ("a(t)") is removed entirely.
$ws.Range("D2").ClearContents()

# E2: "a(t)*u+b(t)*v" -> "u+v"
$ws.Range("E2").Value = "u+v"

# D3 ("b(t)") cleared, keeping its existing style/format.
$ws.Range("D3").ClearContents()

# E3: "c(t)*u+d(t)*v" -> "u+v"
$ws.Range("E3").Value = "u+v"

# D4 ("c(t)") cleared.
$ws.Range("D4").ClearContents()

# D5 ("d(t)") cleared.
$ws.Range("D5").ClearContents()

# Move the active selection to E3, matching the edited view state.
$ws.Range("E3").Select()
